# added key-val pairs for z29 in 32106
# The "key building block" table on Sheet1 is organized into fixed-size
# blocks (20 rows each) per key prefix (m, c1, c2, ... z1, z2, z3). The
# z1 block (rows 291-310) had 17 used rows (z1001-z1017) and 3 free rows
# at the bottom (308-310). Two new key/value pairs (z1018, z1019) were
# added there, which pushed the following z2/z3 blocks down by 10 rows
# (the capacity that was inserted ahead of the z2 block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the z2 (rows 311-320) and z3 (rows 321-332) blocks down by 10 rows
# so the z1 block keeps its own reserved capacity intact.
$ws.Rows("311:320").Insert() | Out-Null

# New pair 1: z1018
$ws.Range("A308").Value = "z1018"
$ws.Range("B308").Value = "직선을 따라 평행이동한 합동인 포물선에 대해 주어진 길이를 알아냅니다."

# New pair 2: z1019
$ws.Range("A309").Value = "z1019"
$ws.Range("B309").Value = "포물선의 정의를 이용해 주어진 선분의 길이에 대한 관계식을 정리해서 요구된 식를 알아냅니다."

# Match the existing column-B formatting (wrapped, s="3") used throughout
# the table by copying it from the row right above the new pairs.
$ws.Range("B307").Copy() | Out-Null
$ws.Range("B308:B309").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection where the author ended up after the edit.
$ws.Range("B312").Select() | Out-Null

Write-Host "Added z1018/z1019 key-value pairs"
